$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision for Friday's movie has ended without a selection.`n"
$ws.Range("D2").Value = "no_decision, "
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding the movie selection for Friday.`n"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("D4").Value = "no_decision, "
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday's showing.`n"
$ws.Range("D5").Value = "no_decision, "
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded: `"Oppenheimer`" will be shown on Friday.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision about Friday's movie cannot be made at this time.`n"
$ws.Range("D7").Value = "no_decision, "
$ws.Range("C8").Value = "MSG: None`n`nMSG: The committee has decided to acquire the rights for `"Oppenheimer`" to be shown on Friday.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for both movies.`n"
$ws.Range("D9").Value = "both_movies, "
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been made to acquire rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for showing on Friday.`n"
$ws.Range("D11").Value = "no_decision, "
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision has been recorded; the movie `"Barbie`" will be shown on Friday.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision to select a movie for Friday has not been made, and the conversation ended without any explicit agreement on a film.`n"
$ws.Range("D14").Value = "no_decision, "
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision-making process ended without a clear choice for Friday's movie, resulting in no decision being made.`n"
$ws.Range("D15").Value = "no_decision, "
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision to acquire a movie for Friday has ended without a resolution.`n"
$ws.Range("D16").Value = "no_decision, "
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday has not been made, so the outcome is recorded as `"no decision.`"`n"
$ws.Range("D17").Value = "no_decision, "
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision about Friday’s movie cannot be made at this time.`n"
$ws.Range("D18").Value = "no_decision, "
$ws.Range("C19").Value = "MSG: None`n`nMSG: The rights for the movie `"Barbie`" have been successfully acquired.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The function has been successfully called, indicating that no decision was made regarding which movie to show on Friday.`n"
$ws.Range("D20").Value = "no_decision, "
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that there was no agreement reached on which movie to show on Friday.`n"
$ws.Range("D21").Value = "no_decision, "
$ws.Range("C22").Value = "MSG: None`n`nMSG: The movie `"Barbie`" has been successfully selected for acquisition.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie was made during the discussion.`n"
$ws.Range("D23").Value = "no_decision, "
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made, so I have called the function for no decision.`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected.`n"
$ws.Range("D25").Value = "no_decision, "
$ws.Range("C26").Value = "MSG: None`n`nMSG: I have recorded the decision as `"no_decision,`" indicating that no movie was explicitly chosen in the meeting.`n"
$ws.Range("D26").Value = "no_decision, "
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired for showing on Friday.`n"
$ws.Range("D29").Value = "both_movies, "
$ws.Range("C30").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D30").Value = "both_movies, "
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Oppenheimer.`"`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision has been logged as `"no decision`" regarding the movie to be shown on Friday.`n"
$ws.Range("D32").Value = "no_decision, "
$ws.Range("C33").Value = "MSG: None`n`nMSG: The decision has been recorded, and there will be no movie selected for Friday at this time.`n"
$ws.Range("D33").Value = "no_decision, "
$ws.Range("C34").Value = "MSG: None`n`nMSG: The decision to acquire rights for both movies has been recorded.`n"
$ws.Range("D34").Value = "both_movies, "
$ws.Range("C35").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision on which movie to show on Friday.`n"
$ws.Range("D35").Value = "no_decision, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The conversation did not lead to a clear decision regarding which movie to show on Friday, so the outcome is classified as no decision made.`n"
$ws.Range("D36").Value = "no_decision, "
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the selection of a movie for Friday. If you need any further assistance, feel free to ask!`n"
$ws.Range("D37").Value = "no_decision, "
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision regarding the movie has been recorded, and no selection was made.`n"
$ws.Range("D38").Value = "no_decision, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie selection for Friday.`n"
$ws.Range("D39").Value = "no_decision, "
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie was selected for Friday.`n"
$ws.Range("D40").Value = "no_decision, "
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision resulted in no movie being selected for Friday.`n"
$ws.Range("D41").Value = "no_decision, "
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`" regarding the movie for Friday's assembly.`n"
$ws.Range("D42").Value = "no_decision, "
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("D43").Value = "no_decision, "
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("D44").Value = "no_decision, "
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded and the function has been called successfully. There is no selected movie for Friday, as the committee did not reach a decision.`n"
$ws.Range("D45").Value = "no_decision, "
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no_decision,`" indicating that no agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Range("D46").Value = "no_decision, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has not been made.`n"
$ws.Range("D47").Value = "no_decision, "
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selection was made during the meeting.`n"
$ws.Range("D48").Value = "no_decision, "
$ws.Range("C49").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D49").Value = "both_movies, "
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("D50").Value = "no_decision, "
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for `"Barbie.`"`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision from the committee is that no movie has been selected for Friday.`n"
$ws.Range("D52").Value = "no_decision, "
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision to acquire the rights for the movie `"Barbie`" has been successfully recorded.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie.`n"
$ws.Range("D54").Value = "no_decision, "
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision process has concluded without an agreement on which movie to show on Friday.`n"
$ws.Range("D57").Value = "no_decision, "
